$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "input_concentrations": the "eq"/"tot" header labels in B1/C1 were
# swapped (column B is now the "tot" column, column C is now the "eq" one).
# ---------------------------------------------------------------------------
$wsConc = $wb.Worksheets.Item("input_concentrations")
$wsConc.Range("B1").Value = "tot"
$wsConc.Range("C1").Value = "eq"

# ---------------------------------------------------------------------------
# Sheet "equilibrium_concentrations": re-computed numeric results (rows 2-4,
# columns B-G). Values below come from the re-run of the upstream Python
# solver referenced in the commit message. Scientific-notation literals are
# routed through [double]"..." since the interpreter's numeric literal
# parser doesn't accept bare "e" exponents.
# ---------------------------------------------------------------------------
$wsEq = $wb.Worksheets.Item("equilibrium_concentrations")

$wsEq.Range("B2").Value = [double]"7.184971860363533e-06"
$wsEq.Range("C2").Value = 0.005
$wsEq.Range("D2").Value = 0.0009918022097660346
$wsEq.Range("E2").Value = [double]"1.012500102643199e-06"
$wsEq.Range("F2").Value = [double]"1.591549489002558e-10"
$wsEq.Range("G2").Value = [double]"6.294627058970815e-11"

$wsEq.Range("B3").Value = [double]"6.691123005003365e-05"
$wsEq.Range("C3").Value = 0.005
$wsEq.Range("D3").Value = 0.0009236320908629021
$wsEq.Range("E3").Value = [double]"9.4290734341462e-06"
$wsEq.Range("F3").Value = [double]"1.380282658894336e-08"
$wsEq.Range("G3").Value = [double]"6.294627058970823e-10"

$wsEq.Range("B4").Value = 0.007544457667507364
$wsEq.Range("C4").Value = 0.005
$wsEq.Range("D4").Value = 0.001041425064918975
$wsEq.Range("E4").Value = 0.001063158535787499
$wsEq.Range("F4").Value = 0.0001754793658933968
$wsEq.Range("G4").Value = [double]"6.294627058970815e-08"

# ---------------------------------------------------------------------------
# Sheet "L_fractions": re-computed numeric results (rows 2-4, columns B-F).
# ---------------------------------------------------------------------------
$wsFrac = $wb.Worksheets.Item("L_fractions")

$wsFrac.Range("B2").Value = 5.143574928422277
$wsFrac.Range("C2").Value = 0.7184971860363533
$wsFrac.Range("D2").Value = 99.18022097660347
$wsFrac.Range("E2").Value = 0.1012500102643199
$wsFrac.Range("F2").Value = [double]"3.183098978005115e-05"

$wsFrac.Range("B3").Value = 4.174500986276123
$wsFrac.Range("C3").Value = 6.691123005003365
$wsFrac.Range("D3").Value = 92.36320908629021
$wsFrac.Range("E3").Value = 0.94290734341462
$wsFrac.Range("F3").Value = 0.002760565317788672

$wsFrac.Range("B4").Value = 2.122371973976437
$wsFrac.Range("C4").Value = 75.44457667507363
$wsFrac.Range("D4").Value = 10.41425064918975
$wsFrac.Range("E4").Value = 10.63158535787499
$wsFrac.Range("F4").Value = 3.509587317867935

# ---------------------------------------------------------------------------
# Sheet "percent_error": re-computed numeric results (rows 2-4, columns B-C).
# ---------------------------------------------------------------------------
$wsErr = $wb.Worksheets.Item("percent_error")

$wsErr.Range("B2").Value = [double]"3.893912102481689e-14"
$wsErr.Range("C2").Value = 0

$wsErr.Range("B3").Value = [double]"2.597748405275269e-16"
$wsErr.Range("C3").Value = 0

$wsErr.Range("B4").Value = [double]"6.29704621779581e-16"
$wsErr.Range("C4").Value = 0
